# Weekly update: a new price-report row is inserted at the top of the
# data block (row 10), pushing the existing rows (10-24) down by one
# (they become rows 11-25). The sheet's used range grows from A1:R24 to
# A1:R25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 10; Excel shifts rows 10:24 down to 11:25
# and carries the row-10 formatting (incl. the date style on column D)
# down with it.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with this week's record.
$ws.Cells.Item(10, 1).Value  = 7
$ws.Cells.Item(10, 2).Value  = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(10, 3).Value  = 'Ñuble'
$ws.Cells.Item(10, 4).Value  = 44791
$ws.Cells.Item(10, 5).Value  = 16
$ws.Cells.Item(10, 6).Value  = 100112043
$ws.Cells.Item(10, 7).Value  = 'Pepino dulce'
$ws.Cells.Item(10, 8).Value  = 'Cultivar IV Región'
$ws.Cells.Item(10, 9).Value  = 'Primera'
$ws.Cells.Item(10, 10).Value = 80
$ws.Cells.Item(10, 11).Value = 17000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 17500
$ws.Cells.Item(10, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(10, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(10, 16).Value = 972
$ws.Cells.Item(10, 17).Value = 18
$ws.Cells.Item(10, 18).Value = 'Hortaliza'
